$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 32) with the latest retrieved job-number stats.
$row = 32
$ws.Range("A$row").Value = 44345.80937171447
$ws.Range("B$row").Value = 74626
$ws.Range("C$row").Value = 62762
$ws.Range("D$row").Value = 3279
$ws.Range("E$row").Value = 2051
$ws.Range("F$row").Value = 1437
$ws.Range("G$row").Value = 19437
$ws.Range("H$row").Value = 1435
$ws.Range("I$row").Value = 824
$ws.Range("J$row").Value = 213

Write-Output "Appended row $row"
